$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 9).Value = "sv"
$ws.Cells.Item(3, 10).Value = "Statement-opinion"
$ws.Cells.Item(8, 9).Value = "sv"
$ws.Cells.Item(8, 10).Value = "Statement-opinion"
$ws.Cells.Item(28, 9).Value = "sd"
$ws.Cells.Item(28, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(47, 9).Value = "ba"
$ws.Cells.Item(47, 10).Value = "Appreciation"
$ws.Cells.Item(64, 9).Value = "aa"
$ws.Cells.Item(64, 10).Value = "Agree/Accept"
$ws.Cells.Item(65, 9).Value = "sd"
$ws.Cells.Item(65, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(76, 9).Value = "ba"
$ws.Cells.Item(76, 10).Value = "Appreciation"
$ws.Cells.Item(97, 9).Value = "sv"
$ws.Cells.Item(97, 10).Value = "Statement-opinion"
$ws.Cells.Item(99, 9).Value = "sd"
$ws.Cells.Item(99, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(106, 9).Value = "sv"
$ws.Cells.Item(106, 10).Value = "Statement-opinion"
$ws.Cells.Item(117, 9).Value = "sv"
$ws.Cells.Item(117, 10).Value = "Statement-opinion"
$ws.Cells.Item(151, 9).Value = "sd"
$ws.Cells.Item(151, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(153, 9).Value = "sv"
$ws.Cells.Item(153, 10).Value = "Statement-opinion"
$ws.Cells.Item(165, 9).Value = "aa"
$ws.Cells.Item(165, 10).Value = "Agree/Accept"
$ws.Cells.Item(251, 9).Value = "sd"
$ws.Cells.Item(251, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(286, 9).Value = "sd"
$ws.Cells.Item(286, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(287, 9).Value = "sd"
$ws.Cells.Item(287, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(301, 9).Value = "sd"
$ws.Cells.Item(301, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(308, 9).Value = "aa"
$ws.Cells.Item(308, 10).Value = "Agree/Accept"
$ws.Cells.Item(314, 9).Value = "sd"
$ws.Cells.Item(314, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(317, 9).Value = "sv"
$ws.Cells.Item(317, 10).Value = "Statement-opinion"
$ws.Cells.Item(319, 9).Value = "sd"
$ws.Cells.Item(319, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(321, 9).Value = "%"
$ws.Cells.Item(321, 10).Value = "Uninterpretable"
$ws.Cells.Item(322, 9).Value = "qy"
$ws.Cells.Item(322, 10).Value = "Yes-No-Question"
$ws.Cells.Item(340, 9).Value = "sd"
$ws.Cells.Item(340, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(343, 9).Value = "%"
$ws.Cells.Item(343, 10).Value = "Uninterpretable"
$ws.Cells.Item(360, 9).Value = "sd"
$ws.Cells.Item(360, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(382, 9).Value = "ba"
$ws.Cells.Item(382, 10).Value = "Appreciation"
$ws.Cells.Item(383, 9).Value = "sd"
$ws.Cells.Item(383, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(385, 9).Value = "ba"
$ws.Cells.Item(385, 10).Value = "Appreciation"
$ws.Cells.Item(389, 9).Value = "sd"
$ws.Cells.Item(389, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(393, 9).Value = "sv"
$ws.Cells.Item(393, 10).Value = "Statement-opinion"
$ws.Cells.Item(409, 9).Value = "sv"
$ws.Cells.Item(409, 10).Value = "Statement-opinion"
$ws.Cells.Item(420, 9).Value = "sv"
$ws.Cells.Item(420, 10).Value = "Statement-opinion"
$ws.Cells.Item(438, 9).Value = "sd"
$ws.Cells.Item(438, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(439, 9).Value = "ba"
$ws.Cells.Item(439, 10).Value = "Appreciation"
$ws.Cells.Item(441, 9).Value = "ba"
$ws.Cells.Item(441, 10).Value = "Appreciation"
$ws.Cells.Item(445, 9).Value = "sv"
$ws.Cells.Item(445, 10).Value = "Statement-opinion"
$ws.Cells.Item(446, 9).Value = "sd"
$ws.Cells.Item(446, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(447, 9).Value = "sd"
$ws.Cells.Item(447, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(449, 9).Value = "sv"
$ws.Cells.Item(449, 10).Value = "Statement-opinion"
$ws.Cells.Item(455, 9).Value = "sd"
$ws.Cells.Item(455, 10).Value = "Statement-non-opinion"
